$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.85916797905456965
$ws.Range("D2").Value = 0.59795320049435663
$ws.Range("A3").Value = 0.87427054050787989
$ws.Range("B3").Value = 0.7653831398886104
$ws.Range("D3").Value = 0.77266934191379866
$ws.Range("E3").Value = 0.79687680462957899
$ws.Range("F4").Value = 0.93068252786538308
$ws.Range("X4").Value = 0.9924352353862762
$ws.Range("AI4").Value = 0.70417178192720542
$ws.Range("F5").Value = 0.83358842423022006
$ws.Range("G5").Value = 0.87281363235480047
$ws.Range("H6").Value = 0.9102277435506021
$ws.Range("AJ6").Value = 0.90355522074314298
$ws.Range("AK7").Value = 0.67257119886871641
$ws.Range("G8").Value = 0.97866711914687587
$ws.Range("Q8").Value = 0.60414407008942383
$ws.Range("AD8").Value = 0.65084582641105104
$ws.Range("G9").Value = 0.85734738327340021
$ws.Range("K9").Value = 0.93995088018433393
$ws.Range("I10").Value = 0.80714763627957953
$ws.Range("K10").Value = 0.99810476251380575
$ws.Range("P10").Value = 0.86351055417943656
$ws.Range("M11").Value = 0.76860424433458108
$ws.Range("Y11").Value = 0.73189978247698817
$ws.Range("L13").Value = 0.66268237009506925
$ws.Range("N13").Value = 0.86943060360412694
$ws.Range("L14").Value = 0.76024636113790378
$ws.Range("P14").Value = 0.97916851611990086
$ws.Range("AO14").Value = 0.59660095349700615
$ws.Range("M15").Value = 0.79740619866094842
$ws.Range("S15").Value = 0.75075214047126415
$ws.Range("I16").Value = 0.95998055980135266
$ws.Range("O17").Value = 0.8516331769158092
$ws.Range("BI17").Value = 0.58580682658871197
$ws.Range("T18").Value = 0.84234830985441833
$ws.Range("BP18").Value = 0.70397469401481683
$ws.Range("Q19").Value = 0.9515325490961547
$ws.Range("T19").Value = 0.99152019916259426
$ws.Range("BE19").Value = 0.93398753472754481
$ws.Range("V20").Value = 0.54961918715341929
$ws.Range("AD20").Value = 0.97133062821538352
$ws.Range("W21").Value = 0.98080479580904201
$ws.Range("BE21").Value = 0.94596865388666784
$ws.Range("W22").Value = 0.68999621051874149
$ws.Range("X23").Value = 0.6624923500901938
$ws.Range("BC23").Value = 0.96721824908902854
$ws.Range("V24").Value = 0.73335349540994343
$ws.Range("W25").Value = 0.82422458598798265
$ws.Range("Z25").Value = 0.82458452082303357
$ws.Range("AA25").Value = 0.82076748907109054
$ws.Range("AL25").Value = 0.68514107706743022
$ws.Range("Z27").Value = 0.98991838428398138
$ws.Range("AB27").Value = 0.79145926175640013
$ws.Range("AK27").Value = 0.95336393081411075
$ws.Range("Z28").Value = 0.60086128536792671
$ws.Range("AD28").Value = 0.65656971836254607
$ws.Range("AA29").Value = 0.82786227591601835
$ws.Range("AB29").Value = 0.60909389816906834
$ws.Range("AN29").Value = 0.78464045205451471
$ws.Range("AC30").Value = 0.80667253313853959
$ws.Range("BB30").Value = 0.78355326604572229
$ws.Range("AC31").Value = 0.72287573736592181
$ws.Range("AD31").Value = 0.75146301339877175
$ws.Range("AF31").Value = 0.68997222814640424
$ws.Range("AG31").Value = 0.87738395355368592
$ws.Range("AD32").Value = 0.51956019283742105
$ws.Range("AG32").Value = 0.9584418073299934
$ws.Range("AH32").Value = 0.89922189901549532
$ws.Range("Z33").Value = 0.65911659504736342
$ws.Range("AH33").Value = 0.8893368428709405
$ws.Range("BI34").Value = 0.97640854317051962
$ws.Range("AJ35").Value = 0.74803454941159653
$ws.Range("AK36").Value = 0.98919309808642519
$ws.Range("AI37").Value = 0.75561732756108224
$ws.Range("X38").Value = 0.92082339003439229
$ws.Range("AJ38").Value = 0.584441109737181
$ws.Range("AN38").Value = 0.89875805935373654
$ws.Range("AM40").Value = 0.76140621673331954
$ws.Range("A41").Value = 0.9381925425548967
$ws.Range("AM41").Value = 0.87759286169800288
$ws.Range("AT41").Value = 0.78643613331122797
$ws.Range("BF41").Value = 0.94342389444674302
$ws.Range("AQ42").Value = 0.92486250370219591
$ws.Range("AR42").Value = 0.68133415572246814
$ws.Range("I43").Value = 0.9296422049591343
$ws.Range("BE43").Value = 0.80912025208515947
$ws.Range("AT44").Value = 0.75622540328321963
$ws.Range("T45").Value = 0.87153891612153722
$ws.Range("AQ45").Value = 0.88704433080782974
$ws.Range("AR45").Value = 0.91250658171368104
$ws.Range("AU45").Value = 0.94068139644406901
$ws.Range("AU46").Value = 0.73991020229752413
$ws.Range("AV46").Value = 0.94388720615482158
$ws.Range("BH47").Value = 0.90032726710260746
$ws.Range("H48").Value = 0.98965260113958575
$ws.Range("AU48").Value = 0.66868295539241007
$ws.Range("AX48").Value = 0.93384825419226103
$ws.Range("AZ48").Value = 0.66398787990041119
$ws.Range("AV49").Value = 0.95907761894310484
$ws.Range("AW50").Value = 0.95739768309839668
$ws.Range("AX51").Value = 0.7507465342240347
$ws.Range("AX52").Value = 0.64201618140456951
$ws.Range("AY52").Value = 0.94887384207234104
$ws.Range("AW53").Value = 0.72700773960604992
$ws.Range("AY53").Value = 0.77972965579252285
$ws.Range("AZ53").Value = 0.97595393051495571
$ws.Range("BP53").Value = 0.90154818270898684
$ws.Range("AP54").Value = 0.62902247387976551
$ws.Range("O55").Value = 0.68236004193086863
$ws.Range("AP56").Value = 0.67737140649596705
$ws.Range("BA56").Value = 0.94333822184137928
$ws.Range("AH57").Value = 0.76014778778111114
$ws.Range("AZ57").Value = 0.77090805822286579
$ws.Range("BD57").Value = 0.83638075482619745
$ws.Range("BF59").Value = 0.84153324377676975
$ws.Range("BG60").Value = 0.77119031795243709
$ws.Range("BI60").Value = 0.76804704924593314
$ws.Range("BJ60").Value = 0.58689379730766822
$ws.Range("AN61").Value = 0.69593549688945422
$ws.Range("BG61").Value = 0.95989334485473565
$ws.Range("BJ61").Value = 0.92607932398541959
$ws.Range("BL62").Value = 0.8353139751012657
$ws.Range("BJ63").Value = 0.7269986875486667
$ws.Range("BL63").Value = 0.99129823858828037
$ws.Range("BM63").Value = 0.90707467494213134
$ws.Range("AJ65").Value = 0.90378134110710662
$ws.Range("BL65").Value = 0.83906384859440153
$ws.Range("AL66").Value = 0.80630305147147863
$ws.Range("BL66").Value = 0.86874694055165969
$ws.Range("A67").Value = 0.65475521206835563
$ws.Range("BM67").Value = 0.83090171063634854
$ws.Range("BN67").Value = 0.90415075880844631
$ws.Range("BP67").Value = 0.70384612906398192
$ws.Range("A68").Value = 0.83344298547462881
$ws.Range("B68").Value = 0.86206740596153419
$ws.Range("J68").Value = 0.90681791050563498
